$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.371.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.608.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.616.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.072.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.291.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.656.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000133"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0742"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.66%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  +5.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.834"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.835"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "275.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.598"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0962"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.951.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.99%  "
